$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2989
$ws1.Range("F3").Value = 6422
$ws1.Range("F6").Value = 543
$ws1.Range("F7").Value = 79
$ws1.Range("F10").Value = 360
$ws1.Range("F12").Value = 7595
$ws1.Range("F13").Value = 365
$ws1.Range("F16").Value = 258
$ws1.Range("F20").Value = 9282
$ws1.Range("F24").Value = 71
$ws1.Range("F27").Value = 123
$ws1.Range("F37").Value = 1490
$ws1.Range("F38").Value = 781
$ws1.Range("F39").Value = 3943
$ws1.Range("F40").Value = 215
$ws1.Range("F41").Value = 45
$ws1.Range("F43").Value = 99
$ws1.Range("F45").Value = 32
$ws1.Range("F48").Value = 39
$ws1.Range("F49").Value = 61

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 53
$ws2.Range("F4").Value = 30
$ws2.Range("F7").Value = 149
$ws2.Range("F15").Value = 6
$ws2.Range("F21").Value = 3

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 53
$ws4.Range("F3").Value = 2989
$ws4.Range("F4").Value = 30
$ws4.Range("F6").Value = 6422
$ws4.Range("F8").Value = 149
$ws4.Range("F10").Value = 543
$ws4.Range("F11").Value = 79
$ws4.Range("F14").Value = 360
$ws4.Range("F17").Value = 7595
$ws4.Range("F18").Value = 365
$ws4.Range("F21").Value = 258
$ws4.Range("F24").Value = 9282
$ws4.Range("F27").Value = 71
$ws4.Range("F29").Value = 123
$ws4.Range("F37").Value = 1490
$ws4.Range("F38").Value = 782
$ws4.Range("F40").Value = 3943
$ws4.Range("F41").Value = 215
$ws4.Range("F42").Value = 45
$ws4.Range("F44").Value = 99
$ws4.Range("F46").Value = 32
$ws4.Range("F48").Value = 39
$ws4.Range("F49").Value = 61
